# Change target field name from "sno" to "seno" in cell D2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "seno"

# Update selection to match the edited cell, as seen in the diff
$ws.Range("D2").Select()
